# Auto-generated edit script: apply updated betting-odds values
# from the 2024-11-11 FlashScore refresh (per row, matching the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("BD2").Value = 151

# Row 6
$ws.Range("G6").Value = 2.4
$ws.Range("I6").Value = 3.2
$ws.Range("J6").Value = 3.25
$ws.Range("AL6").Value = 29
$ws.Range("AS6").Value = 251
$ws.Range("AX6").Value = 19
$ws.Range("BB6").Value = 301

# Row 7
$ws.Range("R7").Value = 1.65

# Row 8
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("Q8").Value = 2.35

# Row 11
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("Z11").Value = 17
$ws.Range("AD11").Value = 7
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 11

# Row 12
$ws.Range("G12").Value = 2.25
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 3.5
$ws.Range("J12").Value = 3.1
$ws.Range("L12").Value = 4.33
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.5
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 2.38
$ws.Range("W12").Value = 6
$ws.Range("X12").Value = 9.5
$ws.Range("Z12").Value = 21
$ws.Range("AA12").Value = 21
$ws.Range("AK12").Value = 41
$ws.Range("AN12").Value = 4
$ws.Range("AO12").Value = 13
$ws.Range("AQ12").Value = 41
$ws.Range("AT12").Value = 2.38
$ws.Range("AX12").Value = 21

# Row 14
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 2.38
$ws.Range("J14").Value = 3.2
$ws.Range("K14").Value = 2.3
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 2.2
$ws.Range("S14").Value = 1.3
$ws.Range("T14").Value = 3.4
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("AA14").Value = 19
$ws.Range("AB14").Value = 23
$ws.Range("AC14").Value = 15
$ws.Range("AD14").Value = 7.5
$ws.Range("AE14").Value = 12
$ws.Range("AG14").Value = 126
$ws.Range("AH14").Value = 11
$ws.Range("AL14").Value = 17
$ws.Range("AT14").Value = 3.4
$ws.Range("AU14").Value = 7
$ws.Range("AY14").Value = 19
$ws.Range("BB14").Value = 101

# Row 15
$ws.Range("G15").Value = 2.4
$ws.Range("I15").Value = 2.7
$ws.Range("Q15").Value = 1.93
$ws.Range("R15").Value = 1.88
$ws.Range("AA15").Value = 21
$ws.Range("AH15").Value = 9
$ws.Range("AI15").Value = 13
$ws.Range("AJ15").Value = 10
$ws.Range("AP15").Value = 23
$ws.Range("AW15").Value = 4.75

# Row 17
$ws.Range("G17").Value = 2.9
$ws.Range("I17").Value = 2.5
$ws.Range("J17").Value = 3.75
$ws.Range("M17").Value = 1.1
$ws.Range("N17").Value = 7
$ws.Range("O17").Value = 1.44
$ws.Range("P17").Value = 2.63
$ws.Range("Q17").Value = 2.4
$ws.Range("R17").Value = 1.53
$ws.Range("S17").Value = 1.53
$ws.Range("T17").Value = 2.38
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 1.73
$ws.Range("Z17").Value = 29
$ws.Range("AC17").Value = 7
$ws.Range("AE17").Value = 17
$ws.Range("AG17").Value = 1000
$ws.Range("AH17").Value = 7
$ws.Range("AI17").Value = 11
$ws.Range("AN17").Value = 4.75
$ws.Range("AO17").Value = 17
$ws.Range("AP17").Value = 29
$ws.Range("AQ17").Value = 51
$ws.Range("AS17").Value = 251
$ws.Range("AT17").Value = 2.38
$ws.Range("AU17").Value = 8.5

# Row 18
$ws.Range("G18").Value = 3.7
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 2.2
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 3
$ws.Range("M18").Value = 1.08
$ws.Range("N18").Value = 8
$ws.Range("O18").Value = 1.4
$ws.Range("P18").Value = 2.75
$ws.Range("Q18").Value = 2.35
$ws.Range("R18").Value = 1.57
$ws.Range("S18").Value = 1.5
$ws.Range("T18").Value = 2.5
$ws.Range("U18").Value = 1.91
$ws.Range("V18").Value = 1.8
$ws.Range("W18").Value = 9
$ws.Range("X18").Value = 17
$ws.Range("Z18").Value = 41
$ws.Range("AA18").Value = 34
$ws.Range("AB18").Value = 41
$ws.Range("AD18").Value = 6
$ws.Range("AE18").Value = 15
$ws.Range("AF18").Value = 51
$ws.Range("AG18").Value = 351
$ws.Range("AH18").Value = 6.5
$ws.Range("AI18").Value = 9.5
$ws.Range("AJ18").Value = 9.5
$ws.Range("AK18").Value = 21
$ws.Range("AL18").Value = 21
$ws.Range("AM18").Value = 34
$ws.Range("AN18").Value = 5.5
$ws.Range("AP18").Value = 34
$ws.Range("AQ18").Value = 67
$ws.Range("AR18").Value = 101
$ws.Range("AS18").Value = 251
$ws.Range("AT18").Value = 2.5
$ws.Range("AU18").Value = 8.5
$ws.Range("AV18").Value = 67
$ws.Range("AW18").Value = 4
$ws.Range("AX18").Value = 13
$ws.Range("AY18").Value = 26
$ws.Range("AZ18").Value = 41
$ws.Range("BA18").Value = 67
$ws.Range("BB18").Value = 201

# Row 19
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 3.2
$ws.Range("I19").Value = 4
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("U19").Value = 1.91
$ws.Range("V19").Value = 1.8
$ws.Range("W19").Value = 6.5
$ws.Range("X19").Value = 9
$ws.Range("Z19").Value = 17
$ws.Range("AA19").Value = 17
$ws.Range("AG19").Value = 351
$ws.Range("AO19").Value = 11
$ws.Range("AV19").Value = 67
$ws.Range("AY19").Value = 34
